$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
  $ws = $wb.Worksheets.Item($sheetName)

  # Row 5 (updated)
  $ws.Range("F5").Value = 3059

  # Row 7 (updated)
  $ws.Range("F7").Value = 2248

  # Row 10 (updated)
  $ws.Range("B10").NumberFormat = "@"
  $ws.Range("B10").Value = "2024-03-24"
  $ws.Range("B10").Style = "Normal"
  $ws.Range("C10").Value = "鹰潭·宅舞联萌·随舞动漫派对（免费活动)"
  $ws.Range("D10").Value = "玉清路与象山路交叉口东南角 鹰潭天虹购物中心"
  $ws.Range("E10").Value = "2024.03.24 14:00-03.24 18:00"
  $ws.Range("F10").Value = 0
  $ws.Range("G10").Value = 22.33
  $ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=82434"
  $ws.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202403/oj5AEi4W1709713367717.jpeg"

  # Row 11 (updated)
  $ws.Range("C11").Value = "南昌·CM01动漫游戏博览会"
  $ws.Range("D11").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
  $ws.Range("E11").Value = "2024.03.30 10:00-03.31 17:00"
  $ws.Range("F11").Value = 1067
  $ws.Range("G11").Value = 55
  $ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=81691"
  $ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202402/9cMJMElF1708938074308.png"

  # Row 12 (updated)
  $ws.Range("B12").NumberFormat = "@"
  $ws.Range("B12").Value = "2024-03-30"
  $ws.Range("B12").Style = "Normal"
  $ws.Range("C12").Value = "鹰潭·原×铁×崩only"
  $ws.Range("D12").Value = "南站路24号 回禾酒店(鹰潭火车站店)"
  $ws.Range("E12").Value = "2024.03.30 10:00-03.30 17:00"
  $ws.Range("F12").Value = 33
  $ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=81097"
  $ws.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg"

  # Row 13 (updated)
  $ws.Range("B13").NumberFormat = "@"
  $ws.Range("B13").Value = "2024-03-31"
  $ws.Range("B13").Style = "Normal"
  $ws.Range("C13").Value = "新余·文旅国漫嘉年华暨BM次元盛典"
  $ws.Range("D13").Value = "五一南路与仙女湖大道交叉口西北 老上海风情街白金汉宫"
  $ws.Range("E13").Value = "2024.03.31 10:00-03.31 17:00"
  $ws.Range("F13").Value = 40
  $ws.Range("G13").Value = 60
  $ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=82208"
  $ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202403/aXc6vPDP1709547191851.jpeg"

  # Row 14 (updated)
  $ws.Range("B14").NumberFormat = "@"
  $ws.Range("B14").Value = "2024-03-31"
  $ws.Range("B14").Style = "Normal"
  $ws.Range("C14").Value = "景德镇·宅舞联萌·随舞动漫派对（免费活动)"
  $ws.Range("D14").Value = "经二路与纬二路交叉路口 景德镇市宝龙广场"
  $ws.Range("E14").Value = "2024.03.31 14:00-03.31 18:00"
  $ws.Range("F14").Value = 0
  $ws.Range("G14").Value = 22.33
  $ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=82437"
  $ws.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202403/zcNNsicc1709714035066.jpeg"

  # Row 15 (updated)
  $ws.Range("B15").NumberFormat = "@"
  $ws.Range("B15").Value = "2024-04-04"
  $ws.Range("B15").Style = "Normal"
  $ws.Range("C15").Value = "南昌·创造力动漫游戏嘉年华1.0"
  $ws.Range("D15").Value = "八一桥街道青山南路118号 蓝海会展中心"
  $ws.Range("E15").Value = "2024.04.04 10:00-04.05 17:00"
  $ws.Range("F15").Value = 35
  $ws.Range("G15").Value = 39.9
  $ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=82419"
  $ws.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202403/dSLjsLyX1709709665435.jpeg"

  # Row 16 (updated)
  $ws.Range("B16").NumberFormat = "@"
  $ws.Range("B16").Value = "2024-04-04"
  $ws.Range("B16").Style = "Normal"
  $ws.Range("C16").Value = "赣州·第三届半夏动漫展"
  $ws.Range("D16").Value = "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
  $ws.Range("E16").Value = "2024.04.04 10:00-04.06 17:00"
  $ws.Range("F16").Value = 265
  $ws.Range("G16").Value = 50
  $ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=82235"
  $ws.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202403/4DWZWYGm1709278879159.jpeg"

  # Row 17 (updated)
  $ws.Range("B17").NumberFormat = "@"
  $ws.Range("B17").Value = "2024-04-04"
  $ws.Range("B17").Style = "Normal"
  $ws.Range("C17").Value = "赣州·赣次元·归来国风动漫节"
  $ws.Range("D17").Value = "客家大道568号文清外国语学校旁 赣州市文清外国语学校国际交流中心"
  $ws.Range("E17").Value = "2024.04.04 10:00-04.04 17:00"
  $ws.Range("F17").Value = 293
  $ws.Range("G17").Value = 40
  $ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=82125"
  $ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202403/fIehikk51709705287036.jpeg"

  # Row 18 (updated)
  $ws.Range("B18").NumberFormat = "@"
  $ws.Range("B18").Value = "2024-04-05"
  $ws.Range("B18").Style = "Normal"
  $ws.Range("C18").Value = "抚州·第七届FZ动漫文化节"
  $ws.Range("D18").Value = "迎宾大道288号 凤凰世纪名都大酒店"
  $ws.Range("E18").Value = "2024.04.05 09:30-04.05 17:00"
  $ws.Range("F18").Value = 6
  $ws.Range("G18").Value = 50
  $ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=82381"
  $ws.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202403/Y725SN0G1709694367526.jpeg"

  # Row 19 (new row)
  $ws.Cells.Item(2,1).Copy() | Out-Null
  $ws.Range("A19").PasteSpecial(-4122) | Out-Null
  $ws.Range("A19").Value = 18
  $ws.Range("B19").NumberFormat = "@"
  $ws.Range("B19").Value = "2024-04-06"
  $ws.Range("B19").Style = "Normal"
  $ws.Range("C19").Value = "萍乡·2024DDS国漫盛典"
  $ws.Range("D19").Value = "凤凰街迎宾路18号 鸿凯大酒店"
  $ws.Range("E19").Value = "2024.04.06 10:00-04.06 17:00"
  $ws.Range("F19").Value = 2
  $ws.Range("G19").Value = 30
  $ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=82413"
  $ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202403/Rbu1xyFr1709707847098.jpeg"

  # Row 20 (new row)
  $ws.Cells.Item(2,1).Copy() | Out-Null
  $ws.Range("A20").PasteSpecial(-4122) | Out-Null
  $ws.Range("A20").Value = 19
  $ws.Range("B20").NumberFormat = "@"
  $ws.Range("B20").Value = "2024-04-13"
  $ws.Range("B20").Style = "Normal"
  $ws.Range("C20").Value = "南昌·原X穹X崩only"
  $ws.Range("D20").Value = "丰和北大道299号 新吉花园酒店"
  $ws.Range("E20").Value = "2024.04.13 10:00-04.13 17:00"
  $ws.Range("F20").Value = 98
  $ws.Range("G20").Value = 65
  $ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=80807"
  $ws.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202402/kfK13XvH1709202705153.jpeg"

  # Row 21 (new row)
  $ws.Cells.Item(2,1).Copy() | Out-Null
  $ws.Range("A21").PasteSpecial(-4122) | Out-Null
  $ws.Range("A21").Value = 20
  $ws.Range("B21").NumberFormat = "@"
  $ws.Range("B21").Value = "2024-04-13"
  $ws.Range("B21").Style = "Normal"
  $ws.Range("C21").Value = "南昌·第二届漫拥动漫嘉年华mini"
  $ws.Range("D21").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
  $ws.Range("E21").Value = "2024.04.13 10:00-04.14 18:00"
  $ws.Range("F21").Value = 45
  $ws.Range("G21").Value = 39.9
  $ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=82210"
  $ws.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202402/KYd0bfk11709203777701.png"

  # Row 22 (new row)
  $ws.Cells.Item(2,1).Copy() | Out-Null
  $ws.Range("A22").PasteSpecial(-4122) | Out-Null
  $ws.Range("A22").Value = 21
  $ws.Range("B22").NumberFormat = "@"
  $ws.Range("B22").Value = "2024-04-20"
  $ws.Range("B22").Style = "Normal"
  $ws.Range("C22").Value = "南昌·DSL国风动漫游戏嘉年华"
  $ws.Range("D22").Value = "沿江北路69号 瑞颐大酒店"
  $ws.Range("E22").Value = "2024.04.20 09:00-04.21 17:00"
  $ws.Range("F22").Value = 73
  $ws.Range("G22").Value = 35
  $ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=82107"
  $ws.Range("I22").Value = "//i0.hdslb.com/bfs/openplatform/202402/QDlumVb41708943318282.jpeg"

}

$excel.CutCopyMode = 0
